$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2241.7576
$ws.Range("I15").Value = 2241.7576
$ws.Range("K15").Value = 6725.2728
$ws.Range("M15").Value = -6556.2728
$ws.Range("H19").Value = 4110.6523
$ws.Range("J19").Value = 3587.5
$ws.Range("L19").Value = 3587.5
$ws.Range("N19").Value = -3937.5
$ws.Range("H33").Value = 111.166664
$ws.Range("I33").Value = 104.25
$ws.Range("K33").Value = 104.25
$ws.Range("M33").Value = 124.75
$ws.Range("H61").Value = 1473
$ws.Range("I61").Value = 1473
$ws.Range("K61").Value = 4419
$ws.Range("M61").Value = -4247
$ws.Range("H106").Value = 3830.4285
$ws.Range("J106").Value = 7004
$ws.Range("L106").Value = 7004
$ws.Range("N106").Value = -8266
$ws.Range("H125").Value = 380
$ws.Range("I125").Value = 380
$ws.Range("K125").Value = 3420
$ws.Range("M125").Value = -960
$ws.Range("H138").Value = 3158.3333
$ws.Range("J138").Value = 2156.6365
$ws.Range("L138").Value = 6469.9095
$ws.Range("N138").Value = -16749.9095
$ws.Range("H141").Value = 1039636.44
$ws.Range("I141").Value = 1476109
$ws.Range("J141").Value = 3014.125
$ws.Range("K141").Value = 4428327
$ws.Range("L141").Value = 9042.375
$ws.Range("M141").Value = -4423147
$ws.Range("N141").Value = -19402.375

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1551839.1
$ws.Range("I2").Value = 2115517
$ws.Range("J2").Value = 1724.75
$ws.Range("K2").Value = 2115517
$ws.Range("L2").Value = 1724.75
$ws.Range("M2").Value = -2115404
$ws.Range("N2").Value = -1950.75
$ws.Range("H32").Value = 3512.3677
$ws.Range("I32").Value = 2992.585
$ws.Range("K32").Value = 2992.585
$ws.Range("M32").Value = -2705.585
$ws.Range("H116").Value = 1551839.1
$ws.Range("I116").Value = 2115517
$ws.Range("J116").Value = 1724.75
$ws.Range("K116").Value = 2115517
$ws.Range("L116").Value = 1724.75
$ws.Range("M116").Value = -2113223
$ws.Range("N116").Value = -6312.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1551839.1
$ws.Range("I3").Value = 2115517
$ws.Range("J3").Value = 1724.75
$ws.Range("K3").Value = 2115517
$ws.Range("L3").Value = 1724.75
$ws.Range("M3").Value = -2115403
$ws.Range("N3").Value = -1952.75
$ws.Range("H76").Value = 43749.5
$ws.Range("J76").Value = 43749.5
$ws.Range("L76").Value = 43749.5
$ws.Range("N76").Value = -44379.5
$ws.Range("H79").Value = 43749.5
$ws.Range("J79").Value = 43749.5
$ws.Range("L79").Value = 43749.5
$ws.Range("N79").Value = -45933.5
$ws.Range("H107").Value = 3432.4375
$ws.Range("I107").Value = 3330.2856
$ws.Range("K107").Value = 3330.2856
$ws.Range("M107").Value = -1410.2856
$ws.Range("H132").Value = 123666
$ws.Range("J132").Value = 123666
$ws.Range("L132").Value = 123666
$ws.Range("N132").Value = -133786
$ws.Range("H134").Value = 7070.121
$ws.Range("J134").Value = 3800
$ws.Range("L134").Value = 11400
$ws.Range("N134").Value = -16470

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3907056.5
$ws.Range("I22").Value = 614.44446
$ws.Range("K22").Value = 614.44446
$ws.Range("M22").Value = -264.44446
$ws.Range("H31").Value = 1487.1765
$ws.Range("I31").Value = 723.75
$ws.Range("J31").Value = 1979.7097
$ws.Range("K31").Value = 723.75
$ws.Range("L31").Value = 1979.7097
$ws.Range("M31").Value = -428.75
$ws.Range("N31").Value = -2569.7097
$ws.Range("H34").Value = 1487.1765
$ws.Range("I34").Value = 723.75
$ws.Range("J34").Value = 1979.7097
$ws.Range("K34").Value = 723.75
$ws.Range("L34").Value = 1979.7097
$ws.Range("M34").Value = -521.75
$ws.Range("N34").Value = -2383.7097
$ws.Range("H62").Value = 6382.7144
$ws.Range("I62").Value = 6912.1665
$ws.Range("J62").Value = 3206
$ws.Range("K62").Value = 6912.1665
$ws.Range("L62").Value = 3206
$ws.Range("M62").Value = -6288.1665
$ws.Range("N62").Value = -4454
$ws.Range("H65").Value = 6382.7144
$ws.Range("I65").Value = 6912.1665
$ws.Range("J65").Value = 3206
$ws.Range("K65").Value = 34560.8325
$ws.Range("L65").Value = 16030
$ws.Range("M65").Value = -31440.8325
$ws.Range("N65").Value = -22270
$ws.Range("H107").Value = 580.4286
$ws.Range("I107").Value = 491.16666
$ws.Range("J107").Value = 699.44446
$ws.Range("K107").Value = 491.16666
$ws.Range("L107").Value = 699.44446
$ws.Range("M107").Value = 1428.83334
$ws.Range("N107").Value = -4539.44446
$ws.Range("H116").Value = 40416.668
$ws.Range("J116").Value = 40416.668
$ws.Range("L116").Value = 40416.668
$ws.Range("N116").Value = -49594.668
$ws.Range("H119").Value = 30416.666
$ws.Range("J119").Value = 30416.666
$ws.Range("L119").Value = 30416.666
$ws.Range("N119").Value = -40092.666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 192.5
$ws.Range("I28").Value = 192.5
$ws.Range("K28").Value = 577.5
$ws.Range("M28").Value = -345.5
$ws.Range("H33").Value = 213.11111
$ws.Range("I33").Value = 79.2
$ws.Range("K33").Value = 475.2
$ws.Range("M33").Value = -192.2
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H68").Value = 2608.4138
$ws.Range("I68").Value = 1496
$ws.Range("J68").Value = 2690.8147
$ws.Range("K68").Value = 4488
$ws.Range("L68").Value = 8072.4441
$ws.Range("M68").Value = -3677
$ws.Range("N68").Value = -9694.444100000001
$ws.Range("H71").Value = 2608.4138
$ws.Range("I71").Value = 1496
$ws.Range("J71").Value = 2690.8147
$ws.Range("K71").Value = 13464
$ws.Range("L71").Value = 24217.3323
$ws.Range("M71").Value = -9408
$ws.Range("N71").Value = -32329.3323
$ws.Range("H81").Value = 28911346
$ws.Range("I81").Value = 3266.6667
$ws.Range("J81").Value = 39751876
$ws.Range("K81").Value = 9800.000100000001
$ws.Range("L81").Value = 119255628
$ws.Range("M81").Value = -8677.000100000001
$ws.Range("N81").Value = -119257874
$ws.Range("H84").Value = 28911346
$ws.Range("I84").Value = 3266.6667
$ws.Range("J84").Value = 39751876
$ws.Range("K84").Value = 29400.0003
$ws.Range("L84").Value = 357766884
$ws.Range("M84").Value = -23784.0003
$ws.Range("N84").Value = -357778116
$ws.Range("H92").Value = 754.55554
$ws.Range("I92").Value = 290.5
$ws.Range("K92").Value = 871.5
$ws.Range("M92").Value = 376.5
$ws.Range("H113").Value = 1438.9231
$ws.Range("J113").Value = 736.4
$ws.Range("L113").Value = 2209.2
$ws.Range("N113").Value = -6549.2
$ws.Range("H131").Value = 6182485
$ws.Range("I131").Value = 166667120
$ws.Range("J131").Value = 9998.833000000001
$ws.Range("K131").Value = 500001360
$ws.Range("L131").Value = 29996.499
$ws.Range("M131").Value = -499996320
$ws.Range("N131").Value = -40076.499

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 635.3333
$ws.Range("I107").Value = 143
$ws.Range("J107").Value = 1250.75
$ws.Range("K107").Value = 143
$ws.Range("L107").Value = 1250.75
$ws.Range("M107").Value = 1777
$ws.Range("N107").Value = -5090.75
$ws.Range("H122").Value = 2504
$ws.Range("J122").Value = 3392.5715
$ws.Range("L122").Value = 10177.7145
$ws.Range("N122").Value = -15077.7145
$ws.Range("H132").Value = 1328716.6
$ws.Range("I132").Value = 1749853.1
$ws.Range("J132").Value = 5144.5713
$ws.Range("K132").Value = 5249559.300000001
$ws.Range("L132").Value = 15433.7139
$ws.Range("M132").Value = -5247029.300000001
$ws.Range("N132").Value = -20493.7139
$ws.Range("H134").Value = 37583.332
$ws.Range("J134").Value = 37583.332
$ws.Range("L134").Value = 112749.996
$ws.Range("N134").Value = -117819.996

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6319.857
$ws.Range("I122").Value = 7024.3335
$ws.Range("J122").Value = 5791.5
$ws.Range("K122").Value = 21073.0005
$ws.Range("L122").Value = 17374.5
$ws.Range("M122").Value = -18623.0005
$ws.Range("N122").Value = -22274.5
$ws.Range("H132").Value = 2949.4888
$ws.Range("I132").Value = 984.3333
$ws.Range("K132").Value = 2952.9999
$ws.Range("M132").Value = -422.9998999999998

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2095.2666
$ws.Range("I81").Value = 1827.4166
$ws.Range("K81").Value = 3654.8332
$ws.Range("M81").Value = -2593.8332
$ws.Range("H84").Value = 2095.2666
$ws.Range("I84").Value = 1827.4166
$ws.Range("K84").Value = 18274.166
$ws.Range("M84").Value = -12970.166
$ws.Range("H126").Value = 5245.273
$ws.Range("I126").Value = 5436.394
$ws.Range("K126").Value = 16309.182
$ws.Range("M126").Value = -13839.182
